$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.696.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.78%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.631.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.77%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.490"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.18%  "
$ws.Range("E8").Value = "  +0.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0619"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.64%  "
$ws.Range("E11").Value = "  +2.42%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.855.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.74%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.625.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.524"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.648.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0732"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "208.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("E22").Value = "  +1.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.21%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.36"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0521"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.07%  "
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.60%  "
$ws.Range("E34").Value = "  +1.46%  "
$ws.Range("E35").Value = "  -0.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.166.51"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0165"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.811"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.14%  "
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.502"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.00%  "
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.766.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0512"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.18%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.61%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.409"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.59%  "
$ws.Range("E51").Value = "  -0.12%  "
